$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update test_number column (A2:A11) from 4 to 2
$ws.Range("A2:A11").Value = 2

# Update the active cell selection to G8
$ws.Range("G8").Select()
